$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range('D2').Value = '51.503.60'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range('D3').Value = '2.984.59'
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '380.43'
$ws.Range('E5').Value = '  +2.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '105.73'
$ws.Range('E6').Value = '  +2.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.542'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range('D9').Value = '0.596'
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '37.44'
$ws.Range('E10').Value = '  +1.56%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range('D13').Value = '3.447.24'
$ws.Range('E13').Value = '  +2.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range('D14').Value = '18.49'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range('D15').Value = '7.58'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '2.972.42'
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range('D17').Value = '0.970'
$ws.Range('E17').Value = '  +2.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range('D18').Value = '51.450.01'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '3.35'
$ws.Range('E19').Value = '  +2.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range('D20').Value = '7.42'
$ws.Range('E20').Value = '  +2.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '12.98'
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0963'
$ws.Range('E22').Value = '  +1.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '69.49'
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '262.06'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range('D25').Value = '2.83'
$ws.Range('E25').Value = '  +5.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '7.59'
$ws.Range('E26').Value = '  +24.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '7.74'
$ws.Range('E27').Value = '  +10.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range('D28').Value = '0.173'
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '0.113'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '25.97'
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '9.93'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range('D33').Value = '35.38'
$ws.Range('E33').Value = '  +1.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '51.30'
$ws.Range('E34').Value = '  +0.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '2.09'
$ws.Range('E35').Value = '  -2.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '0.0448'
$ws.Range('E36').Value = '  +6.26%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range('D38').Value = '3.07'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '17.29'
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '2.61'
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  +2.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range('D43').Value = '125.15'
$ws.Range('E43').Value = '  +4.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '22.22'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '0.292'
$ws.Range('E45').Value = '  +22.42%  '
$ws.Range('E46').Value = '  -0.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '2.37'
$ws.Range('E47').Value = '  +3.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '2.049.22'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '3.26'
$ws.Range('E49').Value = '  +2.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '0.0345'
$ws.Range('E50').Value = '  +10.59%  '
$ws.Range('E51').Value = '  +2.04%  '
